$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the empty "TableCaption"-styled paragraph that sits right before
#    the schedule table.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "TableCaption" -and $p.Range.Text.Trim() -eq "") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Split "Demonstrate the applicable CMS-Required outcomes, state-specific
#    outcomes, or Electronic Visit Verification (EVV) criteria in the
#    production environment." into 5 runs with the revised wording:
#    "Demonstrate the applicable CMS-" / "r" / "equired outcomes" / " and" /
#    " state-specific outcomes in the production environment."
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute(
    "Demonstrate the applicable CMS-Required outcomes, state-specific outcomes, or Electronic Visit Verification (EVV) criteria in the production environment.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2AB72571" w14:textId="77777777" w:rsidR="00A90F1E" w:rsidRPr="00B44D5E" w:rsidRDefault="00A90F1E" w:rsidP="0091485C"><w:pPr><w:pStyle w:val="TableText"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>Demonstrate the applicable CMS-</w:t></w:r><w:r><w:t>r</w:t></w:r><w:r><w:t>equired outcomes</w:t></w:r><w:r><w:t xml:space="preserve"> and</w:t></w:r><w:r><w:t xml:space="preserve"> state-specific outcomes in the production environment.</w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) "Metric/KPI Discussion " -> "Metric Discussion "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Metric/KPI Discussion ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Metric Discussion ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Update the metric/KPI description paragraph wording.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Review the SMC metrics or EVV KPI data (Operational Report Workbook). Discuss any issues or trends that emerged in the metric or KPI data.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Review the metrics data (Operational Report Workbook). Discuss any issues or trends that emerged in the metric data.", 2) | Out-Null
